{"js": "// Add two new paragraphs after the existing \"Third File\" paragraph:\n//   1. an empty paragraph\n//   2. a paragraph with the text \"New Change in Dev Branch\"\n// Both new paragraphs pick up the formatting (Helvetica, 12pt) already in\n// use at the end of the document, matching the existing \"Third File\" run.\nconst body = context.document.body;\n\nconst blankParagraph = body.insertParagraph(\"\", \"End\");\nconst newParagraph = body.insertParagraph(\"New Change in Dev Branch\", \"End\");\n\nawait context.sync();\n", "ps1": "# Add two new paragraphs after the existing \"Third File\" paragraph:\n#   1. an empty paragraph\n#   2. a paragraph with the text \"New Change in Dev Branch\"\n# Both new paragraphs inherit the formatting (Helvetica, 12pt) already used\n# by the last paragraph in the document, matching the \"Third File\" run.\n$d = $word.ActiveDocument\n\n# Insert a blank paragraph right after the current last paragraph.\n$r = $d.Paragraphs.Last.Range\n$r.InsertParagraphAfter()\n\n# Insert another paragraph after that one and give it the new text.\n$r2 = $d.Paragraphs.Last.Range\n$r2.InsertParagraphAfter()\n$r3 = $d.Paragraphs.Last.Range\n$r3.Text = \"New Change in Dev Branch\"\n"}
